$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (column D, may be $null if unchanged),
# new Volume(1h) text (column E).
$updates = @(
    @{ Row = 2;  D = "27.660.92";    E = "  +0.53%  " },
    @{ Row = 3;  D = "1.856.30";     E = "  +0.77%  " },
    @{ Row = 4;  D = "1.031";        E = "  -0.04%  " },
    @{ Row = 5;  D = "322.98";       E = "  +1.20%  " },
    @{ Row = 6;  D = "1.029";        E = "  +0.09%  " },
    @{ Row = 7;  D = "0.4397";       E = "  +0.55%  " },
    @{ Row = 8;  D = "0.3804";       E = "  +1.94%  " },
    @{ Row = 9;  D = "0.07432";      E = "  +0.63%  " },
    @{ Row = 10; D = "0.8849";       E = "  +1.24%  " },
    @{ Row = 11; D = "21.60";        E = "  +0.83%  " },
    @{ Row = 12; D = "1.874.69";     E = "  +1.42%  " },
    @{ Row = 13; D = "5.534";        E = "  +1.03%  " },
    @{ Row = 14; D = "6.736";        E = "  +0.55%  " },
    @{ Row = 15; D = "0.07170";      E = "  +0.13%  " },
    @{ Row = 16; D = "85.37";        E = "  +3.01%  " },
    @{ Row = 17; D = $null;          E = "  +0.15%  " },
    @{ Row = 18; D = "0.000009100";  E = "  +0.94%  " },
    @{ Row = 19; D = $null;          E = "  +0.11%  " },
    @{ Row = 20; D = "15.51";        E = "  +0.46%  " },
    @{ Row = 21; D = "27.701.27";    E = "  +0.60%  " },
    @{ Row = 22; D = "5.302";        E = "  +1.07%  " },
    @{ Row = 23; D = "11.27";        E = "  -0.18%  " },
    @{ Row = 24; D = "2.093.49";     E = "  +0.63%  " },
    @{ Row = 25; D = "2.022";        E = "  +5.55%  " },
    @{ Row = 26; D = "157.97";       E = "  +0.59%  " },
    @{ Row = 27; D = $null;          E = "  +0.76%  " },
    @{ Row = 28; D = "5.370";        E = "  +2.15%  " },
    @{ Row = 29; D = $null;          E = "  +2.87%  " },
    @{ Row = 30; D = "117.92";       E = "  +1.30%  " },
    @{ Row = 31; D = "0.09010";      E = "  -0.69%  " },
    @{ Row = 32; D = "0.7785";       E = "  +2.17%  " },
    @{ Row = 33; D = "1.215";        E = "  +1.09%  " },
    @{ Row = 34; D = "2.998";        E = "  +4.18%  " },
    @{ Row = 35; D = "4.571";        E = "  +1.80%  " },
    @{ Row = 36; D = $null;          E = "  -0.09%  " },
    @{ Row = 37; D = "1.145";        E = "  -0.36%  " },
    @{ Row = 38; D = "0.01979";      E = "  +0.50%  " },
    @{ Row = 39; D = "0.05276";      E = "  +0.30%  " },
    @{ Row = 40; D = "2.857";        E = "  +2.50%  " },
    @{ Row = 41; D = "0.5195";       E = "  +0.39%  " },
    @{ Row = 42; D = "0.1685";       E = "  +1.18%  " },
    @{ Row = 43; D = "6.888";        E = "  +5.12%  " },
    @{ Row = 44; D = "8.860";        E = "  +4.18%  " },
    @{ Row = 45; D = "110.40";       E = "  +1.09%  " },
    @{ Row = 46; D = "10.70";        E = "  +1.28%  " },
    @{ Row = 47; D = "0.06610";      E = "  +4.50%  " },
    @{ Row = 48; D = $null;          E = "  +0.09%  " },
    @{ Row = 49; D = "1.714";        E = "  +1.72%  " },
    @{ Row = 50; D = "0.4720";       E = "  +1.72%  " },
    @{ Row = 51; D = "1.898";        E = "  -0.65%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Column D values look numeric (e.g. "1.031", "85.37") but the
        # source data is plain text, not a real number. Force the cell to
        # text before assigning, then restore its default style so we
        # don't leave a stray number-format style behind.
        $cellD = $ws.Cells.Item($u.Row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
